$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new header columns (I: saved_filename, J: file_path) ---
# Copy the existing header style (bold + border + center/top alignment) from H1
# so the new header cells reuse style index 1 instead of minting a new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "saved_filename"
$ws.Range("J1").Value = "file_path"

# --- Append the two new log rows reflecting the enhanced CSV/logging format ---
$ws.Range("A38").Value = "2025-08-10T16:06:01"
$ws.Range("B38").Value = "admin"
$ws.Range("C38").Value = "Video Upload"
$ws.Range("D38").Value = "Video: 20250810_160601_Khun_A-Standing2.mp4"
$ws.Range("F38").Value = "Khun A-Standing2.mp4"
$ws.Range("G38").Value = 2.84
$ws.Range("H38").Value = "video"
$ws.Range("I38").Value = "20250810_160601_Khun_A-Standing2.mp4"
$ws.Range("J38").Value = "user_uploads/20250810_160601_Khun_A-Standing2.mp4"

$ws.Range("A39").Value = "2025-08-10T16:06:33"
$ws.Range("B39").Value = "admin"
$ws.Range("C39").Value = "Video Upload"
$ws.Range("D39").Value = "Video: 20250810_160633_Pukrirk2.mp4"
$ws.Range("F39").Value = "Pukrirk2.mp4"
$ws.Range("G39").Value = 3.9
$ws.Range("H39").Value = "video"
$ws.Range("I39").Value = "20250810_160633_Pukrirk2.mp4"
$ws.Range("J39").Value = "user_uploads/20250810_160633_Pukrirk2.mp4"
